$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new column before EY, shifting
#     "01-oct." .. "31-oct." (and everything under them) one column to
#     the right, then fill the freshly inserted column with the new
#     "28-dec" header and "-" placeholder values for the 24 data rows.
$ws1 = $wb.Worksheets.Item("Prix Spot")

$ws1.Range("EY:EY").Insert()

$ws1.Range("EY1").Value = "28-dec"
$ws1.Range("EY2:EY25").Value = "-"

# --- Sheet "Gaz": append one new row of data after the existing 182 rows.
#     Column A stores its dates as plain text (e.g. "2025-12-25" on row
#     182), so force text entry for the new date to stop Excel's COM
#     layer from auto-converting the look-alike string into a real date
#     serial value.
$ws2 = $wb.Worksheets.Item("Gaz")

$ws2.Range("A183").NumberFormat = "@"
$ws2.Range("A183").Value = "2025-12-26"
$ws2.Range("A183").ClearFormats()
$ws2.Range("B183").Value = 27.75
